# Fix of the two same subjects
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D3 and D4 were blank placeholders "-", should show the ELT class like C3/C4
$ws.Range("D3").Value = "ELT-2A-Circuitos Elétricos 2"
$ws.Range("D4").Value = "ELT-2A-Circuitos Elétricos 2"

# Row 6: shift the ELT class from C6 to D6 (C6 becomes "-"),
# and move the MCT-2A-Programação class from E6 to B6 (E6 becomes "-")
$ws.Range("B6").Value = "MCT-2A-Programação"
$ws.Range("C6").Value = "-"
$ws.Range("E6").Value = "-"

# Row 7: C7 no longer has the duplicated ELT class
$ws.Range("C7").Value = "-"
